$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '54.152.99'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.262.08'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '495.27'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '128.57'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.524'
$ws.Range('E8').Value = '  -0.92%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('E11').Value = '  +2.73%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.85'
$ws.Range('E12').Value = '  +4.52%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.88'
$ws.Range('E13').Value = '  +5.28%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.662.69'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '54.138.89'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.274.44'
$ws.Range('E17').Value = '  +1.09%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.21'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '302.63'
$ws.Range('E21').Value = '  -1.58%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '60.57'
$ws.Range('E23').Value = '  -2.93%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.998'
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.27'
$ws.Range('E26').Value = '  +3.08%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '172.08'
$ws.Range('E27').Value = '  +2.32%  '
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.08'
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '17.78'
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('E35').Value = '  +2.93%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.20'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.39'
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '124.48'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('E43').Value = '  +1.50%  '
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '240.51'
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.374'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0204'
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('E51').Value = '  -0.40%  '
